# Region XI_HEALTH.xlsx update:
#  - Insert 5 new "No. of Sites ..." summary columns ahead of the existing
#    "Status as of July 4, 2025" column (shifts it from AA to AF).
#  - Re-case two header labels to ALL CAPS.
#  - Drop the stale "-" placeholder values that were sitting in the
#    "Total Physical Target" (I) and "Batch" (L) columns for the first
#    ten data rows (rows 2-11) now that real figures are not yet available.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert five blank columns immediately before the old column AA
# ("Status as of July 4, 2025"); this pushes that column to AF and
# automatically fixes up dimension/dataValidation ranges.
$ws.Range("AA1:AE1").EntireColumn.Insert()

# Re-cased header labels
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# Headers for the newly inserted columns
$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# Clear the leftover "-" placeholders in Total Physical Target / Batch
# for the rows that don't have confirmed values yet.
$ws.Range("I2:I11").ClearContents()
$ws.Range("L2:L11").ClearContents()
